# Update column C ("Förändrad") date values for every data row.
# The "Förändrad" (Changed) date moves from 2023-09-03 (serial 45172)
# to 2023-09-06 (serial 45175).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45175
